# Apply 2025-12-08 violent crime data updates (2025 column) across Citywide Totals,
# By Neighborhood, and individual neighborhood sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 6208
$ws.Range("L3").Value = 6731
$ws.Range("L4").Value = 1670
$ws.Range("L6").Value = 5542
$ws.Range("L7").Value = 20549

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L2").Value = 180
$ws.Range("L4").Value = 72
$ws.Range("L6").Value = 164
$ws.Range("L8").Value = 1356
$ws.Range("L11").Value = 339
$ws.Range("L15").Value = 167
$ws.Range("L19").Value = 555
$ws.Range("L20").Value = 524
$ws.Range("L25").Value = 123
$ws.Range("L27").Value = 176
$ws.Range("L29").Value = 1149
$ws.Range("L31").Value = 204
$ws.Range("L33").Value = 926
$ws.Range("L36").Value = 263
$ws.Range("L42").Value = 655
$ws.Range("L44").Value = 139
$ws.Range("L50").Value = 99
$ws.Range("L52").Value = 438
$ws.Range("L53").Value = 228
$ws.Range("L55").Value = 218
$ws.Range("L56").Value = 21
$ws.Range("L60").Value = 135
$ws.Range("L63").Value = 66
$ws.Range("L64").Value = 129
$ws.Range("L65").Value = 402
$ws.Range("L67").Value = 710
$ws.Range("L74").Value = 20
$ws.Range("L76").Value = 318
$ws.Range("L77").Value = 136
$ws.Range("L79").Value = 566
$ws.Range("L85").Value = 1022
$ws.Range("L89").Value = 283
$ws.Range("L91").Value = 276
$ws.Range("L93").Value = 104
$ws.Range("L99").Value = 361
$ws.Range("L100").Value = 40
$ws.Range("L101").Value = 20549

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("L3").Value = 96
$ws.Range("L7").Value = 339

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("L3").Value = 82
$ws.Range("L6").Value = 81
$ws.Range("L7").Value = 283

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L3").Value = 425
$ws.Range("L6").Value = 211
$ws.Range("L7").Value = 1022

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L2").Value = 137
$ws.Range("L6").Value = 125
$ws.Range("L7").Value = 438

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("L4").Value = 26
$ws.Range("L7").Value = 228

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 411
$ws.Range("L3").Value = 480
$ws.Range("L7").Value = 1356

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L3").Value = 326
$ws.Range("L7").Value = 926

$ws = $wb.Worksheets.Item("New City")
$ws.Range("L6").Value = 98
$ws.Range("L7").Value = 402

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L2").Value = 109
$ws.Range("L7").Value = 361

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("L3").Value = 56
$ws.Range("L7").Value = 204

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L2").Value = 202
$ws.Range("L4").Value = 48
$ws.Range("L7").Value = 710

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L2").Value = 345
$ws.Range("L3").Value = 441
$ws.Range("L4").Value = 62
$ws.Range("L7").Value = 1149

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L6").Value = 151
$ws.Range("L7").Value = 555

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("L4").Value = 11
$ws.Range("L7").Value = 139

$ws = $wb.Worksheets.Item("River North")
$ws.Range("L6").Value = 143
$ws.Range("L7").Value = 318

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("L2").Value = 64
$ws.Range("L7").Value = 164

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L6").Value = 185
$ws.Range("L7").Value = 655

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("L3").Value = 72
$ws.Range("L7").Value = 218

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("L2").Value = 95
$ws.Range("L7").Value = 276

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L4").Value = 40
$ws.Range("L7").Value = 566

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("L3").Value = 38
$ws.Range("L7").Value = 129

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L3").Value = 184
$ws.Range("L7").Value = 524

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("L2").Value = 90
$ws.Range("L7").Value = 263

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("L2").Value = 35
$ws.Range("L7").Value = 104

$ws = $wb.Worksheets.Item("Wrigleyville")
$ws.Range("L6").Value = 25
$ws.Range("L7").Value = 40

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("L2").Value = 41
$ws.Range("L7").Value = 123

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("L2").Value = 63
$ws.Range("L7").Value = 167

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("L6").Value = 27
$ws.Range("L7").Value = 99

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("L2").Value = 59
$ws.Range("L3").Value = 58
$ws.Range("L7").Value = 180

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("L4").Value = 23
$ws.Range("L7").Value = 176

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("L3").Value = 44
$ws.Range("L7").Value = 135

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("L6").Value = 30
$ws.Range("L7").Value = 136

$ws = $wb.Worksheets.Item("Magnificent Mile")
$ws.Range("L6").Value = 10
$ws.Range("L7").Value = 21

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("L2").Value = 27
$ws.Range("L6").Value = 21
$ws.Range("L7").Value = 72

$ws = $wb.Worksheets.Item("Printers Row")
$ws.Range("L6").Value = 8
$ws.Range("L7").Value = 20
